# Generate Report for Handback
#
# Both locale handoffs (zh-cn, de-de) have now come back in sync with
# en-US, so the status report is regenerated:
#   * the shared "Status" text flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it is shown
#     (Overview sheet + each locale sheet),
#   * each locale sheet gets its "Latest Handback DateTime" stamped,
#   * each locale sheet gains the "Latest Target File" (F) / "Latest
#     Handback File" (G) columns, populated + hyperlinked for both
#     tracked source files.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn = $wb.Worksheets.Item(2)
$dede = $wb.Worksheets.Item(3)

$newStatus = "Handed back: in sync with en-US"

# --- Status column everywhere it currently reads "Ready for handoff" ---
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Latest Handback DateTime (column H) per locale ---
$zhcn.Range("H2").Value = "2016-03-22 23:14:01"
$zhcn.Range("H3").Value = "2016-03-22 23:14:01"

$dede.Range("H2").Value = "2016-03-22 23:14:09"
$dede.Range("H3").Value = "2016-03-22 23:14:09"

# --- Latest Target File (F) / Latest Handback File (G) hyperlinked cells ---
$mdName = "4deb35a5-630e-4299-b40e-8f9a90586cb4.md"
$md2Name = "ffffdf981585-01e1-44ae-becf-5b182e7b97a1.md"
$zhcnXlf = "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf"
$dedeXlf = "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/e7f363aaf48fa35c599ab9d8ada0067594646708/e2e/4deb35a5-630e-4299-b40e-8f9a90586cb4.md"
$md2Url = "https://github.com/OpenLocalizationTest/oltest/blob/e7f363aaf48fa35c599ab9d8ada0067594646708/e2e/ffffdf981585-01e1-44ae-becf-5b182e7b97a1.md"
$zhcnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02b07f23aaa38d419af93ba9023d040aa5c23598/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf"
$dedeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/14990f17cbe64192b2773fd6fc146bbe54ec3ff2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf"

# zh-cn sheet: drop every hyperlink then recreate them all (existing +
# new) in row order, so the new F2/G2/F3/G3 links land between the
# existing A/D pairs exactly as the refreshed report lays them out.
$zhcn.Range("A1").Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $mdUrl, "", "", $mdName)
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), $zhcnXlfUrl, "", "", $zhcnXlf)
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $mdUrl, "", "", $mdName)
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), $zhcnXlfUrl, "", "", $zhcnXlf)
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $md2Url, "", "", $md2Name)
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), $zhcnXlfUrl, "", "", $zhcnXlf)
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $mdUrl, "", "", $mdName)
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), $zhcnXlfUrl, "", "", $zhcnXlf)

# de-de sheet: same treatment.
$dede.Range("A1").Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $mdUrl, "", "", $mdName)
$dede.Hyperlinks.Add($dede.Range("D2"), $dedeXlfUrl, "", "", $dedeXlf)
$dede.Hyperlinks.Add($dede.Range("F2"), $mdUrl, "", "", $mdName)
$dede.Hyperlinks.Add($dede.Range("G2"), $dedeXlfUrl, "", "", $dedeXlf)
$dede.Hyperlinks.Add($dede.Range("A3"), $md2Url, "", "", $md2Name)
$dede.Hyperlinks.Add($dede.Range("D3"), $dedeXlfUrl, "", "", $dedeXlf)
$dede.Hyperlinks.Add($dede.Range("F3"), $mdUrl, "", "", $mdName)
$dede.Hyperlinks.Add($dede.Range("G3"), $dedeXlfUrl, "", "", $dedeXlf)
